# Apply price/volume updates from the Aug 26 2023 GitHub Actions crypto data refresh.
# Uses a leading apostrophe on numeric-looking price strings so Excel keeps them as
# text (matching the original inline-string cells) instead of auto-converting to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.152.75"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "1.654.51"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "'217.82"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("E6").Value = "  +1.70%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").Value = "'0.2621"
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("D9").Value = "'0.06318"
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("D10").Value = "'20.38"
$ws.Range("E10").Value = "  -0.30%  "
$ws.Range("D11").Value = "'0.07808"
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("E12").Value = "  +1.15%  "
$ws.Range("D13").Value = "1.645.40"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").Value = "1.883.33"
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").Value = "0.0₅8147"
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("D17").Value = "'65.35"
$ws.Range("E17").Value = "  +1.09%  "
$ws.Range("D18").Value = "26.125.20"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").Value = "'4.598"
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("D21").Value = "'190.61"
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("E22").Value = "  +0.63%  "
$ws.Range("D23").Value = "'5.999"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("D24").Value = "'1.007"
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").Value = "'145.33"
$ws.Range("E25").Value = "  +4.89%  "
$ws.Range("D26").Value = "'0.1224"
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("D27").Value = "'7.200"
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("E28").Value = "  -0.99%  "
$ws.Range("D29").Value = "'1.475"
$ws.Range("E29").Value = "  +5.42%  "
$ws.Range("D30").Value = "'0.05706"
$ws.Range("E30").Value = "  -3.61%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").Value = "'3.546"
$ws.Range("E32").Value = "  +1.57%  "
$ws.Range("D33").Value = "'3.266"
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("D34").Value = "'1.591"
$ws.Range("E34").Value = "  +4.76%  "
$ws.Range("D35").Value = "'2.807"
$ws.Range("E35").Value = "  +2.25%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.422"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'0.9476"
$ws.Range("E37").Value = "  +0.45%  "
$ws.Range("D38").Value = "'0.5715"
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").Value = "'0.01606"
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'5.798"
$ws.Range("E40").Value = "  -0.81%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.8497"
$ws.Range("E41").Value = "  +0.80%  "
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("D43").Value = "'104.08"
$ws.Range("E43").Value = "  +3.62%  "
$ws.Range("D44").Value = "1.038.35"
$ws.Range("E44").Value = "  +3.69%  "
$ws.Range("D45").Value = "1.796.76"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("D46").Value = "'56.71"
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("D47").Value = "0.0₈105"
$ws.Range("E47").Value = "  -1.78%  "
$ws.Range("D48").Value = "'1.004"
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("E49").Value = "  +1.45%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.05154"
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.829"
$ws.Range("E51").Value = "  -0.63%  "
